$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from paragraph 1 ("Left indent") to the
#    empty paragraph right after "What if the third element is in italics?"
#    Word only keeps a single "_GoBack" bookmark, so re-adding it at the new
#    location automatically removes it from the old one.
# ---------------------------------------------------------------------------
$goBackPara = $d.Paragraphs(44)
$d.Bookmarks.Add("_GoBack", $goBackPara.Range)

# ---------------------------------------------------------------------------
# 2) Delete the paragraph that contains only a tab character, just before
#    the "And this is a" numbered-list item.
# ---------------------------------------------------------------------------
$d.Paragraphs(35).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Delete the first of the two empty paragraphs that follow
#    "This is tabbed and italics."
# ---------------------------------------------------------------------------
$d.Paragraphs(29).Range.Delete()

# ---------------------------------------------------------------------------
# 4) Append a bold trailing space, as its own run, to "This is bold text."
#    Splitting a new paragraph off and merging it back in produces a
#    distinct <w:r> instead of getting coalesced into the preceding run.
# ---------------------------------------------------------------------------
$boldPara = $d.Paragraphs(26)
$boldRange = $boldPara.Range
$endOfText = $boldRange.End - 1   # before the paragraph mark
$insertionPoint = $d.Range($endOfText, $endOfText)
$insertionPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs(27)
$newPara.Range.InsertAfter(" ")
$newPara.Range.Font.Bold = 1

$mergeRange = $d.Range($endOfText, $endOfText + 1)
$mergeRange.Delete()

# ---------------------------------------------------------------------------
# 5) Remove the bookmarkStart/bookmarkEnd that used to sit in paragraph 1
#    ("Left indent"). This already happened implicitly in step 1 above
#    (Word relocates rather than duplicates a "_GoBack" bookmark), so
#    nothing further is required here.
# ---------------------------------------------------------------------------

Write-Host "done"
